# Generate Report for Handoff
# Update status/date cells across the three worksheets and resize the
# "date/status" columns that got narrower.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: status + generate-date for both locales ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-10-13 13:50:27"

# --- zh-cn sheet: status + handoff datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-10-13 13:50:16"

# --- de-de sheet: status + handoff datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-10-13 13:50:27"

# --- Column width changes (OOXML stored width 29.9777050018311 -> 17.2159881591797) ---
# Excel's ColumnWidth (character units) round-trips through a pixel-quantised
# formula on save (width = round(chars*6)/6 + 5/6), so the nearest
# representable ColumnWidth for the target stored width is used here.
$newColumnWidth = 16.3333333333333
$wsOverview.Range("E1").ColumnWidth = $newColumnWidth
$wsOverview.Range("F1").ColumnWidth = $newColumnWidth
$wsZhCn.Range("C1").ColumnWidth = $newColumnWidth
$wsDeDe.Range("C1").ColumnWidth = $newColumnWidth
